$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 54 - existing rows 54-81 shift down to 55-82
$ws.Rows.Item(54).Insert()

# Fill the newly inserted row 54 with the new weekly data point
$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44460
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112040
$ws.Range("G54").Value = "Cilantro"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 3300
$ws.Range("K54").Value = 1500
$ws.Range("L54").Value = 2000
$ws.Range("M54").Value = 1750
$ws.Range("N54").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O54").Value = "Provincia del Elquí"
$ws.Range("P54").Value = 1167
$ws.Range("Q54").Value = 1.5
$ws.Range("R54").Value = "Hortaliza"

# Apply the same date-number format used by the other rows in column D
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
